$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 108
$ws.Range("I2").Value = 109.166664
$ws.Range("K2").Value = 109.166664
$ws.Range("M2").Value = 3.833336000000003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 34001.043
$ws.Range("I11").Value = 34001.043
$ws.Range("K11").Value = 34001.043
$ws.Range("M11").Value = -33861.043

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3196.5454
$ws.Range("I62").Value = 2906.2222
$ws.Range("K62").Value = 2906.2222
$ws.Range("M62").Value = -2282.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3196.5454
$ws.Range("I65").Value = 2906.2222
$ws.Range("K65").Value = 14531.111
$ws.Range("M65").Value = -11411.111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 787050.7
$ws.Range("I127").Value = 847323.9
$ws.Range("K127").Value = 2541971.7
$ws.Range("M127").Value = -2537011.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5519.143
$ws.Range("J138").Value = 6015.9
$ws.Range("L138").Value = 18047.7
$ws.Range("N138").Value = -28327.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 464365.3
$ws.Range("I32").Value = 550201.2
$ws.Range("K32").Value = 550201.2
$ws.Range("M32").Value = -549914.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2047.7
$ws.Range("I45").Value = 2153.5
$ws.Range("K45").Value = 2153.5
$ws.Range("M45").Value = -1776.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4076
$ws.Range("I88").Value = 3100
$ws.Range("J88").Value = 7004
$ws.Range("K88").Value = 3100
$ws.Range("L88").Value = 7004
$ws.Range("M88").Value = -2694
$ws.Range("N88").Value = -7816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 4076
$ws.Range("I91").Value = 3100
$ws.Range("J91").Value = 7004
$ws.Range("K91").Value = 3100
$ws.Range("L91").Value = 7004
$ws.Range("M91").Value = -1696
$ws.Range("N91").Value = -9812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2452.6956
$ws.Range("I122").Value = 2366.0588
$ws.Range("J122").Value = 2698.1667
$ws.Range("K122").Value = 7098.176399999999
$ws.Range("L122").Value = 8094.500100000001
$ws.Range("M122").Value = -4648.176399999999
$ws.Range("N122").Value = -12994.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3388.1785
$ws.Range("I132").Value = 2051.35
$ws.Range("K132").Value = 6154.049999999999
$ws.Range("M132").Value = -3624.049999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 463.66666
$ws.Range("I11").Value = 70.5
$ws.Range("K11").Value = 70.5
$ws.Range("M11").Value = 69.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1694.3948
$ws.Range("I94").Value = 1471.7667
$ws.Range("K94").Value = 1471.7667
$ws.Range("M94").Value = -1020.7667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 26885.334
$ws.Range("J103").Value = 26885.334
$ws.Range("L103").Value = 26885.334
$ws.Range("N103").Value = -29229.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5171058
$ws.Range("I134").Value = 4388998.5
$ws.Range("K134").Value = 13166995.5
$ws.Range("M134").Value = -13164460.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 941599.9
$ws.Range("I31").Value = 1794969
$ws.Range("K31").Value = 1794969
$ws.Range("M31").Value = -1794674

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 941599.9
$ws.Range("I34").Value = 1794969
$ws.Range("K34").Value = 1794969
$ws.Range("M34").Value = -1794767

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3346.1
$ws.Range("I132").Value = 3032.8125
$ws.Range("J132").Value = 4599.25
$ws.Range("K132").Value = 9098.4375
$ws.Range("L132").Value = 13797.75
$ws.Range("M132").Value = -6568.4375
$ws.Range("N132").Value = -18857.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6128.174
$ws.Range("I134").Value = 4491.6665
$ws.Range("K134").Value = 13474.9995
$ws.Range("M134").Value = -10939.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3104185.5
$ws.Range("I5").Value = 2101336
$ws.Range("K5").Value = 6304008
$ws.Range("M5").Value = -6303896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 251375.5
$ws.Range("I63").Value = 1829.6666
$ws.Range("K63").Value = 5488.9998
$ws.Range("M63").Value = -4739.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 251375.5
$ws.Range("I66").Value = 1829.6666
$ws.Range("K66").Value = 16466.9994
$ws.Range("M66").Value = -12722.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 13374.875
$ws.Range("I80").Value = 4999
$ws.Range("K80").Value = 14997
$ws.Range("M80").Value = -14061

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 13374.875
$ws.Range("I83").Value = 4999
$ws.Range("K83").Value = 44991
$ws.Range("M83").Value = -40311

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 558093.25
$ws.Range("J129").Value = 2557.7693
$ws.Range("L129").Value = 7673.3079
$ws.Range("N129").Value = -17673.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 3104185.5
$ws.Range("I135").Value = 2101336
$ws.Range("K135").Value = 18912024
$ws.Range("M135").Value = -18909489

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 5840.905
$ws.Range("I20").Value = 4197.2666
$ws.Range("K20").Value = 4197.2666
$ws.Range("M20").Value = -3952.2666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2424
$ws.Range("I102").Value = 2101.2
$ws.Range("K102").Value = 2101.2
$ws.Range("M102").Value = -479.1999999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1011.2222
$ws.Range("I113").Value = 1075.125
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 1075.125
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = 1094.875
$ws.Range("N113").Value = -4840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2483.9473
$ws.Range("J126").Value = 2281.5
$ws.Range("L126").Value = 6844.5
$ws.Range("N126").Value = -11784.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 12569.97
$ws.Range("I132").Value = 7296.7856
$ws.Range("K132").Value = 21890.3568
$ws.Range("M132").Value = -19360.3568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 35000
$ws.Range("J140").Value = 35000
$ws.Range("L140").Value = 35000
$ws.Range("N140").Value = -45360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 73771
$ws.Range("J141").Value = 73771
$ws.Range("L141").Value = 73771
$ws.Range("N141").Value = -84131

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3228.3076
$ws.Range("I7").Value = 2906.9
$ws.Range("K7").Value = 2906.9
$ws.Range("M7").Value = -2794.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1964.5
$ws.Range("I16").Value = 1927.3334
$ws.Range("K16").Value = 1927.3334
$ws.Range("M16").Value = -1757.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3054
$ws.Range("I40").Value = 3019.4
$ws.Range("K40").Value = 3019.4
$ws.Range("M40").Value = -2883.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2791.5833
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3943.4666
$ws.Range("I122").Value = 3653.7144
$ws.Range("K122").Value = 10961.1432
$ws.Range("M122").Value = -8511.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3228.3076
$ws.Range("I126").Value = 2906.9
$ws.Range("K126").Value = 8720.700000000001
$ws.Range("M126").Value = -6250.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 83984.5
$ws.Range("J133").Value = 83984.5
$ws.Range("L133").Value = 83984.5
$ws.Range("N133").Value = -89044.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 14585585
$ws.Range("J136").Value = 16671664
$ws.Range("L136").Value = 50014992
$ws.Range("N136").Value = -50020092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 85145
$ws.Range("J139").Value = 89640
$ws.Range("L139").Value = 89640
$ws.Range("N139").Value = -99920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 1500000
$ws.Range("I21").Value = 1500000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1500000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -1499765
$ws.Range("N21").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 1042500
$ws.Range("I24").Value = 2000000
$ws.Range("K24").Value = 2000000
$ws.Range("M24").Value = -1999770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 85000
$ws.Range("J25").Value = 85000
$ws.Range("L25").Value = 85000
$ws.Range("M25").Value = -85586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 1500000
$ws.Range("I35").Value = 1500000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1500000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1499710
$ws.Range("N35").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2289.4614
$ws.Range("I126").Value = 2348
$ws.Range("K126").Value = 7044
$ws.Range("M126").Value = -4574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 66457
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 66457
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 66457
$ws.Range("M140").Value = $null
$ws.Range("N140").Value = -76817

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 49999
$ws.Range("J141").Value = 49999
$ws.Range("L141").Value = 49999
$ws.Range("N141").Value = -60359
